$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "z1 (1)"

# Update the report title date (15/08/2022 -> 30/06/2022)
$ws.Range("A1").Value = "Report Z1 Measure standard deviations for National; 30/06/2022"

# Update the standard deviation values for 2022 report (30 June)
$ws.Range("C4").Value = 6.1207381871009998
$ws.Range("C5").Value = 5.8760869122460004
$ws.Range("C6").Value = 5.7114030448179998
$ws.Range("C7").Value = 4.9510919988840003
$ws.Range("C8").Value = 8.7174795351449994

# Update the "Generated on" timestamp
$ws.Range("A10").Value = "Generated on 28/09/2023 12:11:18 AEST"
